$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Append 7 new data rows (79-85) to the landscaping log table ---
# All new rows share the same date (2025-05-21 -> serial 45798).
$newDate = 45798

# Give column A on the new rows the same date number format (style) as
# the existing date cells above (row 78) by copying its formatting.
$ws.Range("A79:A85").Value = $newDate
$ws.Cells.Item(78, 1).Copy()
$ws.Range("A79:A85").PasteSpecial(-4122)

# Row 79: Flowering / Large
$ws.Cells.Item(79, 2).Value = "Flowering"
$ws.Cells.Item(79, 3).Value = "Large"
$ws.Cells.Item(79, 4).Value = 51
$ws.Cells.Item(79, 5).Value = 68
$ws.Cells.Item(79, 6).Formula = "=ABS(D79-E79)"
$ws.Cells.Item(79, 7).Value = 3.84
$ws.Cells.Item(79, 8).Value = 2
$ws.Cells.Item(79, 9).Value = "No"
$ws.Cells.Item(79, 10).Value = 2
$ws.Cells.Item(79, 11).Value = "Dark"
$ws.Cells.Item(79, 12).Value = 3
$ws.Cells.Item(79, 13).Value = 0.93
$ws.Cells.Item(79, 14).Value = 62
$ws.Cells.Item(79, 15).Value = 29.6
$ws.Cells.Item(79, 16).Value = 20
$ws.Cells.Item(79, 17).Value = 0.93
$ws.Cells.Item(79, 18).Value = 7.5
$ws.Cells.Item(79, 19).Value = 30
$ws.Cells.Item(79, 20).Value = 37

# Row 80: Nonflowering / Medium
$ws.Cells.Item(80, 2).Value = "Nonflowering"
$ws.Cells.Item(80, 3).Value = "Medium"
$ws.Cells.Item(80, 4).Value = 51
$ws.Cells.Item(80, 5).Value = 68
$ws.Cells.Item(80, 6).Formula = "=ABS(D80-E80)"
$ws.Cells.Item(80, 7).Value = 3.84
$ws.Cells.Item(80, 8).Value = 3
$ws.Cells.Item(80, 9).Value = "No"
$ws.Cells.Item(80, 10).Value = 3
$ws.Cells.Item(80, 11).Value = "Dark"
$ws.Cells.Item(80, 12).Value = 3
$ws.Cells.Item(80, 13).Value = 0.93
$ws.Cells.Item(80, 14).Value = 62
$ws.Cells.Item(80, 15).Value = 29.6
$ws.Cells.Item(80, 16).Value = 20
$ws.Cells.Item(80, 17).Value = 0.93
$ws.Cells.Item(80, 18).Value = 7.5
$ws.Cells.Item(80, 19).Value = 30
$ws.Cells.Item(80, 20).Value = 37

# Row 81: Nonflowering / Small
$ws.Cells.Item(81, 2).Value = "Nonflowering"
$ws.Cells.Item(81, 3).Value = "Small"
$ws.Cells.Item(81, 4).Value = 51
$ws.Cells.Item(81, 5).Value = 68
$ws.Cells.Item(81, 6).Formula = "=ABS(D81-E81)"
$ws.Cells.Item(81, 7).Value = 3.84
$ws.Cells.Item(81, 8).Value = 3.5
$ws.Cells.Item(81, 9).Value = "No"
$ws.Cells.Item(81, 10).Value = 3
$ws.Cells.Item(81, 11).Value = "Dark"
$ws.Cells.Item(81, 12).Value = 3
$ws.Cells.Item(81, 13).Value = 0.93
$ws.Cells.Item(81, 14).Value = 62
$ws.Cells.Item(81, 15).Value = 29.6
$ws.Cells.Item(81, 16).Value = 20
$ws.Cells.Item(81, 17).Value = 0.93
$ws.Cells.Item(81, 18).Value = 7.5
$ws.Cells.Item(81, 19).Value = 30
$ws.Cells.Item(81, 20).Value = 37

# Row 82: Nonflowering / Medium
$ws.Cells.Item(82, 2).Value = "Nonflowering"
$ws.Cells.Item(82, 3).Value = "Medium"
$ws.Cells.Item(82, 4).Value = 51
$ws.Cells.Item(82, 5).Value = 68
$ws.Cells.Item(82, 6).Formula = "=ABS(D82-E82)"
$ws.Cells.Item(82, 7).Value = 3.84
$ws.Cells.Item(82, 8).Value = 4
$ws.Cells.Item(82, 9).Value = "No"
$ws.Cells.Item(82, 10).Value = 3
$ws.Cells.Item(82, 11).Value = "Dark"
$ws.Cells.Item(82, 12).Value = 3
$ws.Cells.Item(82, 13).Value = 0.93
$ws.Cells.Item(82, 14).Value = 62
$ws.Cells.Item(82, 15).Value = 29.6
$ws.Cells.Item(82, 16).Value = 20
$ws.Cells.Item(82, 17).Value = 0.93
$ws.Cells.Item(82, 18).Value = 7.5
$ws.Cells.Item(82, 19).Value = 30
$ws.Cells.Item(82, 20).Value = 37

# Row 83: Nonflowering / Medium
$ws.Cells.Item(83, 2).Value = "Nonflowering"
$ws.Cells.Item(83, 3).Value = "Medium"
$ws.Cells.Item(83, 4).Value = 51
$ws.Cells.Item(83, 5).Value = 68
$ws.Cells.Item(83, 6).Formula = "=ABS(D83-E83)"
$ws.Cells.Item(83, 7).Value = 3.84
$ws.Cells.Item(83, 8).Value = 1
$ws.Cells.Item(83, 9).Value = "No"
$ws.Cells.Item(83, 10).Value = 3
$ws.Cells.Item(83, 11).Value = "Dark"
$ws.Cells.Item(83, 12).Value = 3
$ws.Cells.Item(83, 13).Value = 0.93
$ws.Cells.Item(83, 14).Value = 62
$ws.Cells.Item(83, 15).Value = 29.6
$ws.Cells.Item(83, 16).Value = 20
$ws.Cells.Item(83, 17).Value = 0.93
$ws.Cells.Item(83, 18).Value = 7.5
$ws.Cells.Item(83, 19).Value = 30
$ws.Cells.Item(83, 20).Value = 37

# Row 84: Nonflowering / Large
$ws.Cells.Item(84, 2).Value = "Nonflowering"
$ws.Cells.Item(84, 3).Value = "Large"
$ws.Cells.Item(84, 4).Value = 51
$ws.Cells.Item(84, 5).Value = 68
$ws.Cells.Item(84, 6).Formula = "=ABS(D84-E84)"
$ws.Cells.Item(84, 7).Value = 3.84
$ws.Cells.Item(84, 8).Value = 0.5
$ws.Cells.Item(84, 9).Value = "No"
$ws.Cells.Item(84, 10).Value = 4
$ws.Cells.Item(84, 11).Value = "Dark"
$ws.Cells.Item(84, 12).Value = 3
$ws.Cells.Item(84, 13).Value = 0.93
$ws.Cells.Item(84, 14).Value = 62
$ws.Cells.Item(84, 15).Value = 29.6
$ws.Cells.Item(84, 16).Value = 20
$ws.Cells.Item(84, 17).Value = 0.93
$ws.Cells.Item(84, 18).Value = 7.5
$ws.Cells.Item(84, 19).Value = 30
$ws.Cells.Item(84, 20).Value = 37

# Row 85: Tree / Medium
$ws.Cells.Item(85, 2).Value = "Tree"
$ws.Cells.Item(85, 3).Value = "Medium"
$ws.Cells.Item(85, 4).Value = 51
$ws.Cells.Item(85, 5).Value = 68
$ws.Cells.Item(85, 6).Formula = "=ABS(D85-E85)"
$ws.Cells.Item(85, 7).Value = 3.84
$ws.Cells.Item(85, 8).Value = 2.8
$ws.Cells.Item(85, 9).Value = "No"
$ws.Cells.Item(85, 10).Value = 1
$ws.Cells.Item(85, 11).Value = "Dark"
$ws.Cells.Item(85, 12).Value = 3
$ws.Cells.Item(85, 13).Value = 0.93
$ws.Cells.Item(85, 14).Value = 62
$ws.Cells.Item(85, 15).Value = 29.6
$ws.Cells.Item(85, 16).Value = 20
$ws.Cells.Item(85, 17).Value = 0.93
$ws.Cells.Item(85, 18).Value = 7.5
$ws.Cells.Item(85, 19).Value = 30
$ws.Cells.Item(85, 20).Value = 37

# --- Update the sheet view: scroll/selection moves from the old tail of
# the table (R72:R78) up to the header row's now-last column (U2), and
# the previous frozen/scrolled "topLeftCell" is cleared. ---
[void]$ws.Range("U2").Select()
